$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Cells.Item(17, 8).Value2 = 321379.53
$ws.Cells.Item(17, 10).Value2 = 321379.53
$ws.Cells.Item(17, 12).Value2 = 964138.5900000001
$ws.Cells.Item(17, 14).Value2 = -964474.5900000001
# Row 112
$ws.Cells.Item(112, 8).Value2 = 1385.4
$ws.Cells.Item(112, 9).Value2 = 1512.25
$ws.Cells.Item(112, 10).Value2 = 1365.8846
$ws.Cells.Item(112, 11).Value2 = 4536.75
$ws.Cells.Item(112, 12).Value2 = 4097.6538
$ws.Cells.Item(112, 13).Value2 = -3428.75
$ws.Cells.Item(112, 14).Value2 = -6313.6538
# Row 125
$ws.Cells.Item(125, 8).Value2 = 69142.86
$ws.Cells.Item(125, 10).Value2 = 63000
$ws.Cells.Item(125, 12).Value2 = 567000
$ws.Cells.Item(125, 14).Value2 = -571920
# Row 129
$ws.Cells.Item(129, 8).Value2 = 1372646.8
$ws.Cells.Item(129, 10).Value2 = 2470555.5
$ws.Cells.Item(129, 12).Value2 = 7411666.5
$ws.Cells.Item(129, 14).Value2 = -7421666.5
# Row 132
$ws.Cells.Item(132, 8).Value2 = 3473.2778
$ws.Cells.Item(132, 9).Value2 = 5280.3
$ws.Cells.Item(132, 10).Value2 = 1214.5
$ws.Cells.Item(132, 11).Value2 = 15840.9
$ws.Cells.Item(132, 12).Value2 = 3643.5
$ws.Cells.Item(132, 13).Value2 = -13310.9
$ws.Cells.Item(132, 14).Value2 = -8703.5
# Row 137
$ws.Cells.Item(137, 8).Value2 = 1398.05
$ws.Cells.Item(137, 9).Value2 = 1025.6111
$ws.Cells.Item(137, 10).Value2 = 4750
$ws.Cells.Item(137, 11).Value2 = 3076.8333
$ws.Cells.Item(137, 12).Value2 = 14250
$ws.Cells.Item(137, 13).Value2 = -526.8333000000002
$ws.Cells.Item(137, 14).Value2 = -19350

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value2 = 4922569
$ws.Cells.Item(32, 9).Value2 = 5397819
$ws.Cells.Item(32, 11).Value2 = 5397819
$ws.Cells.Item(32, 13).Value2 = -5397532
# Row 61
$ws.Cells.Item(61, 8).Value2 = 1455.75
$ws.Cells.Item(61, 9).Value2 = 941
$ws.Cells.Item(61, 11).Value2 = 941
$ws.Cells.Item(61, 13).Value2 = -729
# Row 97
$ws.Cells.Item(97, 8).Value2 = 548.3684
$ws.Cells.Item(97, 9).Value2 = 501.11765
$ws.Cells.Item(97, 10).Value2 = 950
$ws.Cells.Item(97, 11).Value2 = 501.11765
$ws.Cells.Item(97, 12).Value2 = 950
$ws.Cells.Item(97, 13).Value2 = -5.117650000000026
$ws.Cells.Item(97, 14).Value2 = -1942
# Row 132
$ws.Cells.Item(132, 8).Value2 = 3057.182
$ws.Cells.Item(132, 9).Value2 = 2329
$ws.Cells.Item(132, 11).Value2 = 6987
$ws.Cells.Item(132, 13).Value2 = -4457
# Row 136
$ws.Cells.Item(136, 8).Value2 = 1455.75
$ws.Cells.Item(136, 9).Value2 = 941
$ws.Cells.Item(136, 11).Value2 = 2823
$ws.Cells.Item(136, 13).Value2 = -273

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Cells.Item(107, 8).Value2 = 10380.066
$ws.Cells.Item(107, 9).Value2 = 1715.4615
$ws.Cells.Item(107, 10).Value2 = 66700
$ws.Cells.Item(107, 11).Value2 = 1715.4615
$ws.Cells.Item(107, 12).Value2 = 66700
$ws.Cells.Item(107, 13).Value2 = 204.5385000000001
$ws.Cells.Item(107, 14).Value2 = -70540

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 50
$ws.Cells.Item(50, 8).Value2 = 9059.125
$ws.Cells.Item(50, 10).Value2 = 9059.125
$ws.Cells.Item(50, 12).Value2 = 9059.125
$ws.Cells.Item(50, 14).Value2 = -10309.125
# Row 51
$ws.Cells.Item(51, 8).Value2 = 8924.5
$ws.Cells.Item(51, 10).Value2 = 8924.5
$ws.Cells.Item(51, 12).Value2 = 8924.5
$ws.Cells.Item(51, 14).Value2 = -10396.5
# Row 61
$ws.Cells.Item(61, 8).Value2 = 8924.5
$ws.Cells.Item(61, 10).Value2 = 8924.5
$ws.Cells.Item(61, 12).Value2 = 8924.5
$ws.Cells.Item(61, 14).Value2 = -9620.5
# Row 99
$ws.Cells.Item(99, 8).Value2 = 3287.5
$ws.Cells.Item(99, 9).Value2 = 4900
$ws.Cells.Item(99, 10).Value2 = 2750
$ws.Cells.Item(99, 11).Value2 = 4900
$ws.Cells.Item(99, 12).Value2 = 2750
$ws.Cells.Item(99, 13).Value2 = -3402
$ws.Cells.Item(99, 14).Value2 = -5746
# Row 126
$ws.Cells.Item(126, 8).Value2 = 3287.5
$ws.Cells.Item(126, 9).Value2 = 4900
$ws.Cells.Item(126, 10).Value2 = 2750
$ws.Cells.Item(126, 11).Value2 = 14700
$ws.Cells.Item(126, 12).Value2 = 8250
$ws.Cells.Item(126, 13).Value2 = -12230
$ws.Cells.Item(126, 14).Value2 = -13190

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Cells.Item(131, 8).Value2 = 7654313.5
$ws.Cells.Item(131, 9).Value2 = 45917828
$ws.Cells.Item(131, 10).Value2 = 1610.6833
$ws.Cells.Item(131, 11).Value2 = 137753484
$ws.Cells.Item(131, 12).Value2 = 4832.0499
$ws.Cells.Item(131, 13).Value2 = -137748444
$ws.Cells.Item(131, 14).Value2 = -14912.0499

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Cells.Item(122, 8).Value2 = 10819821
$ws.Cells.Item(122, 9).Value2 = 7983645.5
$ws.Cells.Item(122, 10).Value2 = 25000700
$ws.Cells.Item(122, 11).Value2 = 23950936.5
$ws.Cells.Item(122, 12).Value2 = 75002100
$ws.Cells.Item(122, 13).Value2 = -23948486.5
$ws.Cells.Item(122, 14).Value2 = -75007000
# Row 132
$ws.Cells.Item(132, 8).Value2 = 10420171
$ws.Cells.Item(132, 9).Value2 = 20837786
$ws.Cells.Item(132, 10).Value2 = 2557.5
$ws.Cells.Item(132, 11).Value2 = 62513358
$ws.Cells.Item(132, 12).Value2 = 7672.5
$ws.Cells.Item(132, 13).Value2 = -62510828
$ws.Cells.Item(132, 14).Value2 = -12732.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Cells.Item(7, 8).Value2 = 79367360
$ws.Cells.Item(7, 9).Value2 = 142859940
$ws.Cells.Item(7, 10).Value2 = 15874788
$ws.Cells.Item(7, 11).Value2 = 142859940
$ws.Cells.Item(7, 12).Value2 = 15874788
$ws.Cells.Item(7, 13).Value2 = -142859828
$ws.Cells.Item(7, 14).Value2 = -15875012
# Row 40
$ws.Cells.Item(40, 8).Value2 = 674940.7
$ws.Cells.Item(40, 9).Value2 = 1264088.8
$ws.Cells.Item(40, 10).Value2 = 1628.5714
$ws.Cells.Item(40, 11).Value2 = 1264088.8
$ws.Cells.Item(40, 12).Value2 = 1628.5714
$ws.Cells.Item(40, 13).Value2 = -1263952.8
$ws.Cells.Item(40, 14).Value2 = -1900.5714
# Row 55
$ws.Cells.Item(55, 8).Value2 = 328.85715
$ws.Cells.Item(55, 9).Value2 = 0
$ws.Cells.Item(55, 10).Value2 = 328.85715
$ws.Cells.Item(55, 11).Value2 = 0
$ws.Cells.Item(55, 12).Value2 = 328.85715
$ws.Cells.Item(55, 13).ClearContents()
$ws.Cells.Item(55, 14).Value2 = -674.85715
# Row 122
$ws.Cells.Item(122, 8).Value2 = 100000
$ws.Cells.Item(122, 9).Value2 = 100000
$ws.Cells.Item(122, 10).Value2 = 0
$ws.Cells.Item(122, 11).Value2 = 300000
$ws.Cells.Item(122, 12).Value2 = 0
$ws.Cells.Item(122, 13).Value2 = -297550
$ws.Cells.Item(122, 14).ClearContents()
# Row 126
$ws.Cells.Item(126, 8).Value2 = 79367360
$ws.Cells.Item(126, 9).Value2 = 142859940
$ws.Cells.Item(126, 10).Value2 = 15874788
$ws.Cells.Item(126, 11).Value2 = 428579820
$ws.Cells.Item(126, 12).Value2 = 47624364
$ws.Cells.Item(126, 13).Value2 = -428577350
$ws.Cells.Item(126, 14).Value2 = -47629304

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Cells.Item(62, 8).Value2 = 8255.556
$ws.Cells.Item(62, 10).Value2 = 10060
$ws.Cells.Item(62, 12).Value2 = 10060
$ws.Cells.Item(62, 14).Value2 = -11308
# Row 65
$ws.Cells.Item(65, 8).Value2 = 8255.556
$ws.Cells.Item(65, 10).Value2 = 10060
$ws.Cells.Item(65, 12).Value2 = 50300
$ws.Cells.Item(65, 14).Value2 = -56540
# Row 122
$ws.Cells.Item(122, 8).Value2 = 1980
$ws.Cells.Item(122, 9).Value2 = 1980
$ws.Cells.Item(122, 10).Value2 = 0
$ws.Cells.Item(122, 11).Value2 = 5940
$ws.Cells.Item(122, 12).Value2 = 0
$ws.Cells.Item(122, 13).Value2 = -3490
$ws.Cells.Item(122, 14).ClearContents()
